$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the "From" rule (row 8) to reflect the jgit commit update.
$ws.Range("E8").Value = "GIT UPDATE"

# Mark E8 as the active/selected cell, matching the saved selection state.
$ws.Range("E8").Select()
